$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 (weekly update), shifting existing rows 12-27 down to 13-28
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the latest weekly price record
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(12, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(12, 4).Value = 44893
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = 100112003
$ws.Cells.Item(12, 7).Value = "Ajo"
$ws.Cells.Item(12, 8).Value = "Chino"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 1400
$ws.Cells.Item(12, 11).Value = 15000
$ws.Cells.Item(12, 12).Value = 16000
$ws.Cells.Item(12, 13).Value = 15571
$ws.Cells.Item(12, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(12, 15).Value = "China"
$ws.Cells.Item(12, 16).Value = 1557
$ws.Cells.Item(12, 17).Value = 10
$ws.Cells.Item(12, 18).Value = "Hortaliza"
